$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Force the new cells to be stored as text so values like dates and
# numeric-looking strings (with leading/trailing spaces) are preserved verbatim.
$newRange = $ws.Range("E13:M14")
$newRange.NumberFormat = "@"

$ws.Range("E13").Value = "09-10-2024"
$ws.Range("F13").Value = "345"
$ws.Range("G13").Value = "P23"
$ws.Range("H13").Value = "696"
$ws.Range("I13").Value = "854"
$ws.Range("J13").Value = "gowtham "
$ws.Range("K13").Value = "test009 "
$ws.Range("L13").Value = "Submitted"
$ws.Range("M13").Value = "Other (Entry manually)"

$ws.Range("E14").Value = "09-10-2024"
$ws.Range("F14").Value = "345"
$ws.Range("G14").Value = "P23"
$ws.Range("H14").Value = "777888 "
$ws.Range("I14").Value = "854 "
$ws.Range("J14").Value = "tester  "
$ws.Range("K14").Value = "90"
$ws.Range("L14").Value = "Submitted"
$ws.Range("M14").Value = "Material schedule No/Service Reason"
